$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original column A values (14 rows) before mutating anything,
# so column B can be populated as a verbatim copy of the pre-edit column A.
$orig = @()
for ($r = 1; $r -le 14; $r++) {
    $orig += $ws.Cells.Item($r, 1).Value2
}

# Add a second column (B) that duplicates the original column A contents.
for ($r = 1; $r -le 14; $r++) {
    $ws.Cells.Item($r, 2).Value = $orig[$r - 1]
}

# Relabel A1 (was "Functions.NavNavigate.page") with the new text, and make
# both new header cells (A1, B1) bold.
$ws.Range("A1").Value = "Nav.Navigate.page"
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Font.Bold = $true

# Column width tweaks: widen column A a bit, and give the new column C the
# width column A used to have; column B stays at the sheet's default width.
$ws.Columns.Item(1).ColumnWidth = 26
$ws.Columns.Item(3).ColumnWidth = 24.1

# Move the active selection to A2 (was A6).
$ws.Range("A2").Select()
